# daily auto push: 2026-02-15 18:53 UTC
#
# Inserts one new data row into Sheet1 at row 817 (pushing the existing
# rows 817-858 down to 818-859) with the values:
#   A817 = "2026/02/16"  (text, matches the "yyyy/mm/dd" text style of the
#                          rest of column A - NOT an Excel date serial)
#   B817 = "月"
#   C817 = 0
#   D817 = 201
#
# This grows the sheet's used range from A1:D858 to A1:D859.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 817; everything at/after 817
# (the 2026/12/29 ... 2027/01/05 block) shifts down by one row.
$ws.Rows.Item(817).Insert()

# Column A throughout the sheet stores dates as plain text (e.g.
# "2026/02/16"), not as real Excel dates. Assigning a date-shaped string
# directly would get auto-recognised and converted to a date serial, so
# force the cell to text format first, then restore "General" so no
# stray formatting lingers on the cell.
$ws.Range("A817").NumberFormat = "@"
$ws.Range("A817").Value = "2026/02/16"
$ws.Range("A817").NumberFormat = "General"

$ws.Range("B817").Value = "月"
$ws.Range("C817").Value = 0
$ws.Range("D817").Value = 201
